$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (translated to human-readable Spanish labels)
$ws.Range("A1").Value = "Número de edificios"
$ws.Range("B1").Value = "Viviendas en el edificio"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Comarca código"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Municipio código"
$ws.Range("H1").Value = "Provincia nombre"
$ws.Range("I1").Value = "Municipio nombre"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:numero-de-edificios"
$ws.Range("B2").Value = "iaest-dimension:viviendas-en-el-edificio"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("I2").Value = "sdmx-dimension:refArea"

# Row 3 - role (medida/dim)
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "dim"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "null"
$ws.Range("H3").Value = "dim"
$ws.Range("I3").Value = "dim"

# Row 4 - type / URI references
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "skos:Concept"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "null"
$ws.Range("H4").Value = "URI-Provincia"
$ws.Range("I4").Value = "URI-Municipio"

# Row 5 - mapping file moved from column F to column B (copy preserves the
# existing cell style so the destination keeps s="1" like the source had)
$ws.Range("F5").Copy($ws.Range("B5"))
$ws.Range("F5").Clear()
